$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.4100175597556301
$ws.Cells.Item(2, 4).Value = 0.03408206458156826
$ws.Cells.Item(2, 5).Value = 0.1733727720169114
$ws.Cells.Item(2, 6).Value = 0.6606048217810709
$ws.Cells.Item(2, 7).Value = 0.5010273403105145
$ws.Cells.Item(2, 8).Value = 0.6520397677357721
$ws.Cells.Item(2, 11).Value = 1.853689686677853
$ws.Cells.Item(2, 12).Value = 0.1514483768234953
$ws.Cells.Item(2, 15).Value = 2.270256797398801

$ws.Cells.Item(3, 3).Value = 0.4008109535730568
$ws.Cells.Item(3, 4).Value = 0.03158881245200718
$ws.Cells.Item(3, 5).Value = 0.1690780971998862
$ws.Cells.Item(3, 6).Value = 0.6660892639240785
$ws.Cells.Item(3, 7).Value = 0.5081027607491535
$ws.Cells.Item(3, 8).Value = 0.6606238608951216
$ws.Cells.Item(3, 11).Value = 1.633340012826352
$ws.Cells.Item(3, 12).Value = 0.1473393954464797
$ws.Cells.Item(3, 15).Value = 2.302900786523963

$ws.Cells.Item(4, 3).Value = 0.3953787099069643
$ws.Cells.Item(4, 4).Value = 0.03005079516928788
$ws.Cells.Item(4, 5).Value = 0.166537162879461
$ws.Cells.Item(4, 6).Value = 0.6700626017487394
$ws.Cells.Item(4, 7).Value = 0.5129985167585005
$ws.Cells.Item(4, 8).Value = 0.6663237581009298
$ws.Cells.Item(4, 11).Value = 1.497567283811691
$ws.Cells.Item(4, 12).Value = 0.144902541993801
$ws.Cells.Item(4, 15).Value = 2.325000857630101

$ws.Cells.Item(5, 3).Value = 0.3932205692147193
$ws.Cells.Item(5, 4).Value = 0.02942228815950187
$ws.Cells.Item(5, 5).Value = 0.1655258654182639
$ws.Cells.Item(5, 6).Value = 0.6718337479717746
$ws.Cells.Item(5, 7).Value = 0.5151317428898778
$ws.Cells.Item(5, 8).Value = 0.6687543177877586
$ws.Cells.Item(5, 11).Value = 1.442122628280174
$ws.Cells.Item(5, 12).Value = 0.1439311453029859
$ws.Cells.Item(5, 15).Value = 2.334522627130667

$ws.Cells.Item(6, 3).Value = 0.39286556799442
$ws.Cells.Item(6, 4).Value = 0.02931782055455301
$ws.Cells.Item(6, 5).Value = 0.1653593994080396
$ws.Cells.Item(6, 6).Value = 0.6721370136841642
$ws.Cells.Item(6, 7).Value = 0.5154942938640161
$ws.Cells.Item(6, 8).Value = 0.6691644177960399
$ws.Cells.Item(6, 11).Value = 1.432909182485218
$ws.Cells.Item(6, 12).Value = 0.1437711525743879
$ws.Cells.Item(6, 15).Value = 2.336134824677501

$ws.Cells.Item(7, 3).Value = 0.3953493794978158
$ws.Cells.Item(7, 4).Value = 0.03004232593946199
$ws.Cells.Item(7, 5).Value = 0.1665234263552797
$ws.Cells.Item(7, 6).Value = 0.6700858731848882
$ws.Cells.Item(7, 7).Value = 0.513026727419664
$ws.Cells.Item(7, 8).Value = 0.6663561011513792
$ws.Cells.Item(7, 11).Value = 1.496820003058247
$ws.Cells.Item(7, 12).Value = 0.1448893537797673
$ws.Cells.Item(7, 15).Value = 2.325127184803009

$ws.Cells.Item(8, 3).Value = 0.4067973570493848
$ws.Cells.Item(8, 4).Value = 0.0332239017205751
$ws.Cells.Item(8, 5).Value = 0.1718720450062463
$ws.Cells.Item(8, 6).Value = 0.6623698855319091
$ws.Cells.Item(8, 7).Value = 0.5033522021724508
$ws.Cells.Item(8, 8).Value = 0.6549104001717012
$ws.Cells.Item(8, 11).Value = 1.777814357124555
$ws.Cells.Item(8, 12).Value = 0.1500137303992872
$ws.Cells.Item(8, 15).Value = 2.281084820413781

$ws.Cells.Item(9, 3).Value = 0.4309967050731132
$ws.Cells.Item(9, 4).Value = 0.03940453814650624
$ws.Cells.Item(9, 5).Value = 0.1831228155443654
$ws.Cells.Item(9, 6).Value = 0.6520633353394985
$ws.Cells.Item(9, 7).Value = 0.4887779859588619
$ws.Cells.Item(9, 8).Value = 0.6358768046682712
$ws.Cells.Item(9, 11).Value = 2.324917802778145
$ws.Cells.Item(9, 12).Value = 0.1607465263643206
$ws.Cells.Item(9, 15).Value = 2.211094800523966

$ws.Cells.Item(10, 3).Value = 0.4498441367013868
$ws.Cells.Item(10, 4).Value = 0.04390802423466766
$ws.Cells.Item(10, 5).Value = 0.1918549804939076
$ws.Cells.Item(10, 6).Value = 0.6474566402747328
$ws.Cells.Item(10, 7).Value = 0.4807819077317106
$ws.Cells.Item(10, 8).Value = 0.6239801422826758
$ws.Cells.Item(10, 11).Value = 2.724334652824837
$ws.Cells.Item(10, 12).Value = 0.1690511796421958
$ws.Cells.Item(10, 15).Value = 2.169738627529227

$ws.Cells.Item(11, 3).Value = 0.4586507327794322
$ws.Cells.Item(11, 4).Value = 0.04594830252867155
$ws.Cells.Item(11, 5).Value = 0.195929124009723
$ws.Cells.Item(11, 6).Value = 0.6460101370822358
$ws.Cells.Item(11, 7).Value = 0.4777396789714743
$ws.Cells.Item(11, 8).Value = 0.6190229301029717
$ws.Cells.Item(11, 11).Value = 2.90546024572717
$ws.Cells.Item(11, 12).Value = 0.172920793680774
$ws.Cells.Item(11, 15).Value = 2.15312778676801

$ws.Cells.Item(12, 3).Value = 0.462019023754209
$ws.Cells.Item(12, 4).Value = 0.04671965981414417
$ws.Cells.Item(12, 5).Value = 0.1974865503754799
$ws.Cells.Item(12, 6).Value = 0.6455561390012079
$ws.Cells.Item(12, 7).Value = 0.4766737939097538
$ws.Cells.Item(12, 8).Value = 0.6172112829848118
$ws.Cells.Item(12, 11).Value = 2.973962521182443
$ws.Cells.Item(12, 12).Value = 0.1743993385658911
$ws.Cells.Item(12, 15).Value = 2.147155821881029

$ws.Cells.Item(13, 3).Value = 0.4612921163136434
$ws.Cells.Item(13, 4).Value = 0.04655359063494302
$ws.Cells.Item(13, 5).Value = 0.1971504798909365
$ws.Cells.Item(13, 6).Value = 0.6456497388733027
$ws.Cells.Item(13, 7).Value = 0.4768995114794166
$ws.Cells.Item(13, 8).Value = 0.6175985360022622
$ws.Cells.Item(13, 11).Value = 2.959213219843718
$ws.Cells.Item(13, 12).Value = 0.1740803200088976
$ws.Cells.Item(13, 15).Value = 2.148427815282332

$ws.Cells.Item(14, 3).Value = 0.458927174535944
$ws.Cells.Item(14, 4).Value = 0.04601178791121185
$ws.Cells.Item(14, 5).Value = 0.1960569611163407
$ws.Cells.Item(14, 6).Value = 0.6459709048680011
$ws.Cells.Item(14, 7).Value = 0.4776502584400646
$ws.Cells.Item(14, 8).Value = 0.6188725703094491
$ws.Cells.Item(14, 11).Value = 2.911097716156178
$ws.Cells.Item(14, 12).Value = 0.17304216974685
$ws.Cells.Item(14, 15).Value = 2.152630083765899

$ws.Cells.Item(15, 3).Value = 0.4574829306632182
$ws.Cells.Item(15, 4).Value = 0.04567975369525357
$ws.Cells.Item(15, 5).Value = 0.1953890551951289
$ws.Cells.Item(15, 6).Value = 0.6461798512540042
$ws.Cells.Item(15, 7).Value = 0.4781213469585595
$ws.Cells.Item(15, 8).Value = 0.6196614934488167
$ws.Cells.Item(15, 11).Value = 2.881614258159914
$ws.Cells.Item(15, 12).Value = 0.1724079927723068
$ws.Cells.Item(15, 15).Value = 2.1552455784323

$ws.Cells.Item(16, 3).Value = 0.4492732850307277
$ws.Cells.Item(16, 4).Value = 0.04377451460810278
$ws.Cells.Item(16, 5).Value = 0.1915907742235774
$ws.Cells.Item(16, 6).Value = 0.647564270328381
$ws.Cells.Item(16, 7).Value = 0.480992748913792
$ws.Cells.Item(16, 8).Value = 0.6243132703655334
$ws.Cells.Item(16, 11).Value = 2.712485839813212
$ws.Cells.Item(16, 12).Value = 0.1688001388507843
$ws.Cells.Item(16, 15).Value = 2.170868630200232

$ws.Cells.Item(17, 3).Value = 0.4442965224603483
$ws.Cells.Item(17, 4).Value = 0.04260353231352809
$ws.Cells.Item(17, 5).Value = 0.1892867295708669
$ws.Cells.Item(17, 6).Value = 0.6485801180217905
$ws.Cells.Item(17, 7).Value = 0.482907103460299
$ws.Cells.Item(17, 8).Value = 0.6272835462695809
$ws.Cells.Item(17, 11).Value = 2.608581941619491
$ws.Cells.Item(17, 12).Value = 0.1666103510035981
$ws.Cells.Item(17, 15).Value = 2.181018003771015

$ws.Cells.Item(18, 3).Value = 0.4414559360908186
$ws.Cells.Item(18, 4).Value = 0.04192922824915257
$ws.Cells.Item(18, 5).Value = 0.1879710899279914
$ws.Cells.Item(18, 6).Value = 0.649225478594353
$ws.Cells.Item(18, 7).Value = 0.4840641825240084
$ws.Cells.Item(18, 8).Value = 0.6290347533410667
$ws.Cells.Item(18, 11).Value = 2.548765641102079
$ws.Cells.Item(18, 12).Value = 0.165359481796429
$ws.Cells.Item(18, 15).Value = 2.187062854701452

$ws.Cells.Item(19, 3).Value = 0.4404979280389227
$ws.Cells.Item(19, 4).Value = 0.04170078691419832
$ws.Cells.Item(19, 5).Value = 0.1875272835710291
$ws.Cells.Item(19, 6).Value = 0.6494544628970758
$ws.Cells.Item(19, 7).Value = 0.4844655497359227
$ws.Cells.Item(19, 8).Value = 0.6296350249500193
$ws.Cells.Item(19, 11).Value = 2.528503796796144
$ws.Cells.Item(19, 12).Value = 0.1649374424838328
$ws.Cells.Item(19, 15).Value = 2.189145077941362

$ws.Cells.Item(20, 3).Value = 0.4448240395122696
$ws.Cells.Item(20, 4).Value = 0.04272826706859689
$ws.Cells.Item(20, 5).Value = 0.1895310070085756
$ws.Cells.Item(20, 6).Value = 0.6484656555385371
$ws.Cells.Item(20, 7).Value = 0.4826975181129427
$ws.Cells.Item(20, 8).Value = 0.6269629262585426
$ws.Cells.Item(20, 11).Value = 2.61964825735248
$ws.Cells.Item(20, 12).Value = 0.1668425635182018
$ws.Cells.Item(20, 15).Value = 2.179916130254327

$ws.Cells.Item(21, 3).Value = 0.4596209085596854
$ws.Cells.Item(21, 4).Value = 0.04617096284810884
$ws.Cells.Item(21, 5).Value = 0.1963777569125256
$ws.Cells.Item(21, 6).Value = 0.6458740226322632
$ws.Cells.Item(21, 7).Value = 0.4774274036602577
$ws.Cells.Item(21, 8).Value = 0.6184965754571152
$ws.Cells.Item(21, 11).Value = 2.925232772500124
$ws.Cells.Item(21, 12).Value = 0.1733467411763172
$ws.Cells.Item(21, 15).Value = 2.151387127748819

$ws.Cells.Item(22, 3).Value = 0.4694863153303004
$ws.Cells.Item(22, 4).Value = 0.04841364569891482
$ws.Cells.Item(22, 5).Value = 0.2009378139636269
$ws.Cells.Item(22, 6).Value = 0.6447269471535861
$ws.Cells.Item(22, 7).Value = 0.4744854276555941
$ws.Cells.Item(22, 8).Value = 0.6133454285873512
$ws.Cells.Item(22, 11).Value = 3.124446533126786
$ws.Cells.Item(22, 12).Value = 0.1776745756921798
$ws.Cells.Item(22, 15).Value = 2.134597188163866

$ws.Cells.Item(23, 3).Value = 0.4642031554209041
$ws.Cells.Item(23, 4).Value = 0.04721736853389302
$ws.Cells.Item(23, 5).Value = 0.1984962222329045
$ws.Cells.Item(23, 6).Value = 0.645289000728539
$ws.Cells.Item(23, 7).Value = 0.4760094700879307
$ws.Cells.Item(23, 8).Value = 0.6160596757009245
$ws.Cells.Item(23, 11).Value = 3.018169744278168
$ws.Cells.Item(23, 12).Value = 0.1753576829386247
$ws.Cells.Item(23, 15).Value = 2.143388025482352

$ws.Cells.Item(24, 3).Value = 0.4445854848196404
$ws.Cells.Item(24, 4).Value = 0.04267187787974791
$ws.Cells.Item(24, 5).Value = 0.1894205411882552
$ws.Cells.Item(24, 6).Value = 0.6485172129990673
$ws.Cells.Item(24, 7).Value = 0.4827920957221608
$ws.Cells.Item(24, 8).Value = 0.6271077429013161
$ws.Cells.Item(24, 11).Value = 2.614645426670904
$ws.Cells.Item(24, 12).Value = 0.1667375551078436
$ws.Cells.Item(24, 15).Value = 2.18041363373473

$ws.Cells.Item(25, 3).Value = 0.4242626606873046
$ws.Cells.Item(25, 4).Value = 0.0377389565087185
$ws.Cells.Item(25, 5).Value = 0.1799974267339906
$ws.Cells.Item(25, 6).Value = 0.6543325083138214
$ws.Cells.Item(25, 7).Value = 0.492246728967551
$ws.Cells.Item(25, 8).Value = 0.6406598933237433
$ws.Cells.Item(25, 11).Value = 2.177347470515599
$ws.Cells.Item(25, 12).Value = 0.1577695626390465
$ws.Cells.Item(25, 15).Value = 2.228267236292709
